$d = $word.ActiveDocument

# 1) "Week Number:" -> add " 7" right after the colon.
$rng = $d.Content
$found = $rng.Find.Execute("Week Number:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter(" 7")
}

# 2) "Milestone Number: " -> add "1" right after the existing trailing space.
$rng = $d.Content
$found = $rng.Find.Execute("Milestone Number: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.InsertAfter("1")
}

# 3) "Project Name:" -> add " Product Ranking Website" (bold, matching the
#    existing run formatting) right after the colon.
$rng = $d.Content
$found = $rng.Find.Execute("Project Name:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Collapse(0)
    $rng.Font.Bold = 1
    $rng.InsertAfter(" Product Ranking Website")
}
